$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.157.35"
$ws.Range("E2").Value = "  +1.22%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.216.09"
$ws.Range("E3").Value = "  -0.56%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.83"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "60.81"
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.14%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0899"
$ws.Range("E10").Value = "  +2.28%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.542.51"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.45"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.13"
$ws.Range("E14").Value = "  +2.60%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.234.62"
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.122.82"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0928"
$ws.Range("E19").Value = "  +4.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.18"
$ws.Range("E20").Value = "  +2.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.02"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "243.57"
$ws.Range("E22").Value = "  -1.80%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.41"
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.38"
$ws.Range("E25").Value = "  +5.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.60"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "169.05"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "20.35"
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.46"
$ws.Range("E30").Value = "  +2.62%  "
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  -1.06%  "
$ws.Range("E33").Value = "  -1.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.61"
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0648"
$ws.Range("E35").Value = "  +3.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.56"
$ws.Range("E36").Value = "  -3.98%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.31"
$ws.Range("E37").Value = "  -3.96%  "
$ws.Range("E38").Value = "  -1.76%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0249"
$ws.Range("E39").Value = "  +5.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.000228"
$ws.Range("E41").Value = "  -4.19%  "
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0956"
$ws.Range("E43").Value = "  -2.85%  "
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.20"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "96.88"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("E46").Value = "  -10.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.456.33"
$ws.Range("E47").Value = "  -0.97%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "15.98"
$ws.Range("E48").Value = "  -2.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.07"
$ws.Range("E49").Value = "  -1.13%  "
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.20"
$ws.Range("E51").Value = "  +1.27%  "
